# #327 Ajout des profils d'acces
# - Update the "Date" metadata value
# - Swap the order of the two "Mapping" columns (RIM Mapping / Spécification
#   métier) in the Elements sheet, including header text, data values and
#   column widths

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet: bump the recorded generation Date
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# ---------------------------------------------------------------------
# 2) Elements sheet: swap columns AK (37) and AL (38)
#    Before: AK = "Mapping: RIM Mapping", AL = "Mapping: Spécification métier..."
#    After : AK = "Mapping: Spécification métier...", AL = "Mapping: RIM Mapping"
#    (the underlying data for each row travels with its column header)
# ---------------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

$lastRow = 6
for ($r = 1; $r -le $lastRow; $r++) {
    $akCell = $elements.Cells.Item($r, 37)
    $alCell = $elements.Cells.Item($r, 38)
    $akValue = $akCell.Value2
    $alValue = $alCell.Value2
    # Skip cells that would not actually change (e.g. both sides already
    # blank) so we don't needlessly disturb untouched rows.
    if ($akValue -ne $alValue) {
        $akCell.Value2 = $alValue
        $alCell.Value2 = $akValue
    }
}

# ---------------------------------------------------------------------
# 3) Swap the column widths that go along with columns AK and AL so the
#    now-wide "Spécification métier" column (AK) and the now-narrow
#    "RIM Mapping" column (AL) keep displaying correctly.
#    (ColumnWidth is expressed in characters; the original stored OOXML
#    <col> widths were 24.98046875 (AK) and 90.640625 (AL); they trade
#    places. The saved <col> width = ColumnWidth + 5/6, so we compensate
#    for that fixed offset when setting the new widths.)
# ---------------------------------------------------------------------
$colAK = $elements.Columns.Item(37)
$colAL = $elements.Columns.Item(38)

$colAK.ColumnWidth = 89.83333333333333
$colAL.ColumnWidth = 24.166666666666668
